$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the sheet: insert a new row for "naive_bayes" (new row 6), pushing
#     the old row 6 ("logistic_regression") down to row 7. Copy formats from
#     row 6 into row 7 first so the new row inherits identical cell styling
#     (e.g. column A's bordered/bold/centered style) before values change.
$ws.Range("A6:F6").Copy($ws.Range("A7:F7"))

# The "accuracy" column (D) holds numeric-looking values that are stored as
# TEXT in this sheet (t="inlineStr" in the original XML), so force the
# column to Text format for the duration of the write, then drop back to the
# default "Normal" style (no explicit number format) to match the sheet's
# existing un-styled data cells.
$accCells = $ws.Range("D2:D7")
$accCells.NumberFormat = "@"

# --- Row 2 ---
$ws.Range("D2").Value = "0.6004999999999999"
$ws.Range("F2").Value = 28.96503901481628

# --- Row 3 ---
$ws.Range("D3").Value = "0.6004999999999999"
$ws.Range("F3").Value = 2.488220691680908

# --- Row 4 ---
$ws.Range("D4").Value = "0.1975"
$ws.Range("F4").Value = 4.691967964172363

# --- Row 5: experiment renamed (naive_bayes -> adaboost) ---
$ws.Range("C5").Value = "scargc_JITC_adaboost_UNSW"
$ws.Range("D5").Value = "0.6004999999999999"
$ws.Range("F5").Value = 31.99111533164978

# --- Row 6: new naive_bayes row ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "JITC"
$ws.Range("C6").Value = "scargc_JITC_naive_bayes_UNSW"
$ws.Range("D6").Value = "0.1975"
$ws.Range("E6").Value = "Only one class found"
$ws.Range("F6").Value = 1.328789710998535

# --- Row 7: old row 6 (logistic_regression), refreshed metrics ---
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "JITC"
$ws.Range("C7").Value = "scargc_JITC_logistic_regression_UNSW"
$ws.Range("D7").Value = "0.48575"
$ws.Range("E7").Value = "Only one class found"
$ws.Range("F7").Value = 6.465306520462036

# Drop the transient Text number format back to the default style so column D
# cells end up with no explicit style, matching the source workbook.
$accCells.Style = "Normal"
